# Each of the 5 section-id markers in this document is split across three
# runs: "<id>" (Courier New, gold), the bare id text (default formatting),
# and "</id>" (Courier New, gold). Collapse each triple into a single run
# "<id>p066r_N</id>" that keeps the first run's (Courier New) formatting.
# A self Find&Replace over the full visible text of each triplet merges the
# matched runs into one, inheriting the formatting of the first run.
$d = $word.ActiveDocument

for ($i = 1; $i -le 5; $i++) {
    $needle = "<id>p066r_$i</id>"
    $found = $d.Content.Find.Execute($needle, $true, $false, $false, $false, $false,
                             $true, 1, $false, $needle, 2)
    Write-Host ("p066r_$i replaced: " + $found)
}
